$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-08 Thursday" "2024-08-09 Friday"

Replace-Text "940×8=7520" "694×3=2082"
Replace-Text "399×2=798" "102×4=408"
Replace-Text "958×2=1916" "577×8=4616"
Replace-Text "770×3=2310" "446×4=1784"
Replace-Text "512×6=3072" "276×5=1380"

Replace-Text "686×4=2744" "528×9=4752"
Replace-Text "454×9=4086" "435×2=870"
Replace-Text "573×9=5157" "799×2=1598"
Replace-Text "649×3=1947" "599×2=1198"
Replace-Text "613×9=5517" "691×3=2073"

Replace-Text "814×9=7326" "211×2=422"
Replace-Text "520×4=2080" "380×9=3420"
Replace-Text "926×2=1852" "693×3=2079"
Replace-Text "633×3=1899" "115×2=230"
Replace-Text "713×8=5704" "350×4=1400"

Replace-Text "502×9=4518" "373×2=746"
Replace-Text "120×6=720" "190×3=570"
Replace-Text "452×5=2260" "544×9=4896"
Replace-Text "133×8=1064" "673×6=4038"
Replace-Text "239×3=717" "885×5=4425"

Replace-Text "648×2=1296" "869×2=1738"
Replace-Text "571×7=3997" "463×3=1389"
Replace-Text "477×9=4293" "803×5=4015"
Replace-Text "769×5=3845" "793×8=6344"
Replace-Text "145×3=435" "776×6=4656"
